$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E15").Value = 121
$ws.Range("F15").Value = 60
$ws.Range("H15").Value = 60

$ws.Range("E17").Value = 78
$ws.Range("F17").Value = 35
$ws.Range("H17").Value = 35

$ws.Range("E18").Value = 75

$ws.Range("E29").Value = 11
$ws.Range("F29").Value = 7
$ws.Range("H29").Value = 7

$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 0
$ws.Range("H31").Value = 0

$ws.Range("E33").Value = 21

$ws.Range("E37").Value = 29
$ws.Range("F37").Value = 15
$ws.Range("H37").Value = 15

$ws.Range("E38").Value = 43

$ws.Range("E46").Value = 18

$ws.Range("E70").Value = 22

$ws.Range("E76").Value = 31

$ws.Range("E77").Value = 31
